# "se subio las nuevas clases" -- new attendance ("p" = presente) marks were
# added for two extra class dates (columns K and L) on the Hoja1 attendance
# grid, for every student row (3-14) except row 12 which only gained the
# second date (L12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Rows 3-11 and 13-14 get both K and L marked; row 12 only gets L.
$fullRows = @(3,4,5,6,7,8,9,10,11,13,14)
foreach ($r in $fullRows) {
    $ws.Range("K$r").Value = "p"
    $ws.Range("L$r").Value = "p"
}
$ws.Range("L12").Value = "p"

# View-state bookkeeping to mirror the author's saved selection/zoom.
$ws.Range("L8").Select()
$excel.ActiveWindow.Zoom = 80

$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Activate()
$excel.ActiveWindow.Zoom = 80

$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Activate()
$excel.ActiveWindow.Zoom = 80

# Restore Hoja1 as the active/selected sheet.
$ws.Activate()
$ws.Range("L8").Select()

try {
    $excel.ActiveWindow.TabRatio = 0.821
} catch {
}
